$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row 4 (Oct 24 2020 vs Sunrisers Hyderabad) entry entirely.
$ws.Rows.Item(4).Delete()

# Remove the row 2 (Oct 30 2020 vs Rajasthan Royals) entry entirely.
# This shifts the former row 3 (Nov 1 2020 vs Chennai Super Kings) up to row 2.
$ws.Rows.Item(2).Delete()
